# Apply edit: add "TestSummary" worksheet with a couple of summary cells,
# and clear the stray empty D15 cell on the "TestResults" sheet.

$wb = $excel.ActiveWorkbook

# --- Clear the leftover empty cell D15 on TestResults sheet ---
$wsResults = $wb.Worksheets.Item("TestResults")
$wsResults.Range("D15").Value = ""

# --- Add the new TestSummary sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsSummary = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsSummary.Name = "TestSummary"

$wsSummary.Range("A1").Value = "Test Executed On"
$wsSummary.Range("B1").Value = "2021-07-03 13:29:51.983579"
$wsSummary.Range("A2").Value = "Number of Test Cases"
$wsSummary.Range("B2").Formula = "=(COUNTA(TestResults!A:A) - 1)"

# Keep the originally active sheet selected (unrelated to the added content)
$wb.Worksheets.Item(1).Activate()
